# fhir ig initial setup
$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the generation timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-08-20T08:30:34+05:45"

# --- Include #0 sheet: refresh the SNOMED CT concept list ---
$inc = $wb.Worksheets.Item("Include #0")

$inc.Cells.Item(2, 1).Value = "308335008"
$inc.Cells.Item(2, 2).Value = "Patient encounter procedure (procedure)"

$inc.Cells.Item(3, 1).Value = "185318008"
$inc.Cells.Item(3, 2).Value = "Third party encounter (procedure)"

$inc.Cells.Item(4, 1).Value = "390906007"
$inc.Cells.Item(4, 2).Value = "Follow-up encounter (procedure)"

# The fourth old concept row is gone; deleting it shifts the blank
# separator row and the "System URI" row up by one (rows 6/7 -> 5/6).
$inc.Rows.Item(5).Delete()
